$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 77 (existing rows 77..97 shift down to 79..99).
$ws.Range("A77:R78").EntireRow.Insert()

# New row 77: Cebollín, Primera, Región de Ñuble, 10-Mar-2023
$ws.Cells.Item(77, 1).Value = 11
$ws.Cells.Item(77, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(77, 3).Value = "Bíobío"
$ws.Cells.Item(77, 4).Value = 44995
$ws.Cells.Item(77, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(77, 5).Value = 8
$ws.Cells.Item(77, 6).Value = 100112037
$ws.Cells.Item(77, 7).Value = "Cebollín"
$ws.Cells.Item(77, 8).Value = "Sin especificar"
$ws.Cells.Item(77, 9).Value = "Primera"
$ws.Cells.Item(77, 10).Value = 200
$ws.Cells.Item(77, 11).Value = 700
$ws.Cells.Item(77, 12).Value = 800
$ws.Cells.Item(77, 13).Value = 750
$ws.Cells.Item(77, 14).Value = "`$/paquete 6 unidades"
$ws.Cells.Item(77, 15).Value = "Región de Ñuble"
$ws.Cells.Item(77, 16).Value = 125
$ws.Cells.Item(77, 17).Value = 6
$ws.Cells.Item(77, 18).Value = "Hortaliza"

# New row 78: Cebollín, Segunda, Región de Ñuble, 10-Mar-2023
$ws.Cells.Item(78, 1).Value = 11
$ws.Cells.Item(78, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(78, 3).Value = "Bíobío"
$ws.Cells.Item(78, 4).Value = 44995
$ws.Cells.Item(78, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(78, 5).Value = 8
$ws.Cells.Item(78, 6).Value = 100112037
$ws.Cells.Item(78, 7).Value = "Cebollín"
$ws.Cells.Item(78, 8).Value = "Sin especificar"
$ws.Cells.Item(78, 9).Value = "Segunda"
$ws.Cells.Item(78, 10).Value = 100
$ws.Cells.Item(78, 11).Value = 600
$ws.Cells.Item(78, 12).Value = 600
$ws.Cells.Item(78, 13).Value = 600
$ws.Cells.Item(78, 14).Value = "`$/paquete 6 unidades"
$ws.Cells.Item(78, 15).Value = "Región de Ñuble"
$ws.Cells.Item(78, 16).Value = 100
$ws.Cells.Item(78, 17).Value = 6
$ws.Cells.Item(78, 18).Value = "Hortaliza"
